$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 87, shifting existing rows 87-144 down to 88-145.
$ws.Rows.Item(87).Insert()

# Fill the new row 87 with its data. Columns A, B, C, E, F, G, H, I, N, O, Q, R
# repeat the same "Apio" template used throughout this sheet.
$ws.Range("A87").Value = 5
$ws.Range("B87").Value = "Macroferia Regional de Talca"
$ws.Range("C87").Value = "Maule"
$ws.Range("D87").Value = 44651
$ws.Range("E87").Value = 7
$ws.Range("F87").Value = 100112017
$ws.Range("G87").Value = "Apio"
$ws.Range("H87").Value = "Americana (o)"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 600
$ws.Range("K87").Value = 7000
$ws.Range("L87").Value = 7000
$ws.Range("M87").Value = 7000
$ws.Range("N87").Value = "`$/docena de matas"
$ws.Range("O87").Value = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P87").Value = 1167
$ws.Range("Q87").Value = 6
$ws.Range("R87").Value = "Hortaliza"
